$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for PQ, ramp, limit trafo, variability function of time
$ws.Range("H1").Value = "pf_ess"
$ws.Range("I1").Value = "ramp"

# Add new values for row 2
$ws.Range("H2").Value = 0.9
$ws.Range("I2").Value = 50

# Update the active selection to reflect where the user left off editing
$ws.Range("H3").Select()
